$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = 2
$ws.Range("F6").Value = 1
$ws.Range("F9").Value = -3
$ws.Range("F10").Value = -1
$ws.Range("F11").Value = 1
$ws.Range("F13").Value = 1
$ws.Range("F14").Value = -1
$ws.Range("F16").Value = -4
$ws.Range("F20").Value = 4
$ws.Range("F21").Value = -1
$ws.Range("F23").Value = 2
